$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.348248600959778
$ws.Range("B1").Value = 3.043877363204956
$ws.Range("C1").Value = 5.132182121276855
$ws.Range("D1").Value = 2.171936273574829
$ws.Range("E1").Value = 1.021070599555969
